$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The previous run only analyzed one file (resume.txt); re-running the
# analysis tool over the full input folder produced results for two files,
# so row 2 gets corrected/updated and a new row 3 is appended with the
# second file's results.
$ws.Range("A2").Value = "Two-Timescale Gradient Descent Ascent Algorithms for.txt"
$ws.Range("B2").Value = 203
$ws.Range("C2").Value = 30
$ws.Range("D2").Value = 0
$ws.Range("E2").Value = 20425

$ws.Range("A3").Value = "A-computational-analysis-of-transcribed-speech-of-people_2025_Computer-Speec.txt"
$ws.Range("B3").Value = 545
$ws.Range("C3").Value = 85
$ws.Range("D3").Value = 0
$ws.Range("E3").Value = 20501
